$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update Critical Minutes values (C3, C4) from 175 -> 143
$ws.Range("C3").Value = 143
$ws.Range("C4").Value = 143

# Populate Driver Vintage for the 22.150.3.1 driver (E12), previously blank.
# Assign as a string literal via a formula then convert to a static value so
# Excel stores it as literal text "2022-08-29" instead of auto-converting the
# date-like text into a date serial number, and so the cell keeps its
# original style (no new number-format style gets introduced).
$ws.Range("E12").Formula = "=""2022-08-29"""
$ws.Range("E12").Copy()
$ws.Range("E12").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false
